$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9831453561782837
$ws.Range("B1").Value = 1.895299792289734
$ws.Range("C1").Value = 5.257884502410889
$ws.Range("D1").Value = 2.229696035385132
$ws.Range("E1").Value = 1.315937161445618
